$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 4124.2
$ws.Range("I19").Value = 600
$ws.Range("J19").Value = 4515.778
$ws.Range("K19").Value = 600
$ws.Range("L19").Value = 4515.778
$ws.Range("M19").Value = -425
$ws.Range("N19").Value = -4865.778
# Row 40
$ws.Range("H40").Value = 3890.3667
$ws.Range("I40").Value = 1189
$ws.Range("J40").Value = 4430.64
$ws.Range("K40").Value = 1189
$ws.Range("L40").Value = 4430.64
$ws.Range("M40").Value = -1014
$ws.Range("N40").Value = -4780.64
# Row 41
$ws.Range("H41").Value = 795.4815
$ws.Range("I41").Value = 709.3333
$ws.Range("K41").Value = 709.3333
$ws.Range("M41").Value = -269.3333
# Row 46
$ws.Range("H46").Value = 747.8333
$ws.Range("I46").Value = 684.5
$ws.Range("J46").Value = 874.5
$ws.Range("K46").Value = 2053.5
$ws.Range("L46").Value = 2623.5
$ws.Range("M46").Value = -1934.5
$ws.Range("N46").Value = -2861.5
# Row 60
$ws.Range("H60").Value = 747.8333
$ws.Range("I60").Value = 684.5
$ws.Range("J60").Value = 874.5
$ws.Range("K60").Value = 2053.5
$ws.Range("L60").Value = 2623.5
$ws.Range("M60").Value = -1569.5
$ws.Range("N60").Value = -3591.5
# Row 64
$ws.Range("H64").Value = 6413.75
$ws.Range("I64").Value = 3400
$ws.Range("K64").Value = 3400
$ws.Range("M64").Value = -3152
# Row 67
$ws.Range("H67").Value = 6413.75
$ws.Range("I67").Value = 3400
$ws.Range("K67").Value = 3400
$ws.Range("M67").Value = -2542
# Row 82
$ws.Range("H82").Value = 1698.8
$ws.Range("J82").Value = 2000
$ws.Range("L82").Value = 6000
$ws.Range("N82").Value = -6812
# Row 85
$ws.Range("H85").Value = 1698.8
$ws.Range("J85").Value = 2000
$ws.Range("L85").Value = 6000
$ws.Range("N85").Value = -8808
# Row 116
$ws.Range("H116").Value = 5608.5557
$ws.Range("I116").Value = 4897.3335
$ws.Range("K116").Value = 4897.3335
$ws.Range("M116").Value = -1455.3335
# Row 132
$ws.Range("H132").Value = 17243860
$ws.Range("I132").Value = 17859624
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 53578872
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -53576342
$ws.Range("N132").Value = -12560
# Row 135
$ws.Range("H135").Value = 2183.5
$ws.Range("I135").Value = 1157.9231
$ws.Range("K135").Value = 10421.3079
$ws.Range("M135").Value = -7886.3079
# Row 138
$ws.Range("H138").Value = 2897.3186
$ws.Range("I138").Value = 1927.5769
$ws.Range("J138").Value = 3285.2153
$ws.Range("K138").Value = 5782.7307
$ws.Range("L138").Value = 9855.6459
$ws.Range("M138").Value = -642.7307000000001
$ws.Range("N138").Value = -20135.6459

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2637.3845
$ws.Range("I61").Value = 2427.3044
$ws.Range("J61").Value = 4248
$ws.Range("K61").Value = 2427.3044
$ws.Range("L61").Value = 4248
$ws.Range("M61").Value = -2215.3044
$ws.Range("N61").Value = -4672
# Row 136
$ws.Range("H136").Value = 2637.3845
$ws.Range("I136").Value = 2427.3044
$ws.Range("J136").Value = 4248
$ws.Range("K136").Value = 7281.9132
$ws.Range("L136").Value = 12744
$ws.Range("M136").Value = -4731.9132
$ws.Range("N136").Value = -17844

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 257.6
$ws.Range("I7").Value = 147
$ws.Range("K7").Value = 147
$ws.Range("M7").Value = -34
# Row 31
$ws.Range("H31").Value = 24293.83
$ws.Range("I31").Value = 1515.3478
$ws.Range("K31").Value = 1515.3478
$ws.Range("M31").Value = -1220.3478
# Row 34
$ws.Range("H34").Value = 24293.83
$ws.Range("I34").Value = 1515.3478
$ws.Range("K34").Value = 1515.3478
$ws.Range("M34").Value = -1313.3478
# Row 50
$ws.Range("H50").Value = 8000
$ws.Range("J50").Value = 8000
$ws.Range("L50").Value = 8000
$ws.Range("N50").Value = -9250
# Row 99
$ws.Range("H99").Value = 4874.25
$ws.Range("I99").Value = 5500
$ws.Range("J99").Value = 4665.6665
$ws.Range("K99").Value = 5500
$ws.Range("L99").Value = 4665.6665
$ws.Range("M99").Value = -4002
$ws.Range("N99").Value = -7661.6665
# Row 126
$ws.Range("H126").Value = 4874.25
$ws.Range("I126").Value = 5500
$ws.Range("J126").Value = 4665.6665
$ws.Range("K126").Value = 16500
$ws.Range("L126").Value = 13996.9995
$ws.Range("M126").Value = -14030
$ws.Range("N126").Value = -18936.9995

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 13389469
$ws.Range("I4").Value = 17462384
$ws.Range("K4").Value = 52387152
$ws.Range("M4").Value = -52387040
# Row 44
$ws.Range("H44").Value = 143638.28
$ws.Range("J44").Value = 333999.66
$ws.Range("L44").Value = 1001998.98
$ws.Range("N44").Value = -1002794.98
# Row 52
$ws.Range("H52").Value = 1039.3334
$ws.Range("J52").Value = 1039.3334
$ws.Range("L52").Value = 3118.0002
$ws.Range("N52").Value = -3650.0002
# Row 68
$ws.Range("H68").Value = 878.2222
$ws.Range("I68").Value = 684
$ws.Range("J68").Value = 975.3333
$ws.Range("K68").Value = 2052
$ws.Range("L68").Value = 2925.9999
$ws.Range("M68").Value = -1241
$ws.Range("N68").Value = -4547.9999
# Row 69
$ws.Range("H69").Value = 3834.3333
$ws.Range("I69").Value = 3504.5
$ws.Range("J69").Value = 3999.25
$ws.Range("K69").Value = 10513.5
$ws.Range("L69").Value = 11997.75
$ws.Range("M69").Value = -9702.5
$ws.Range("N69").Value = -13619.75
# Row 71
$ws.Range("H71").Value = 878.2222
$ws.Range("I71").Value = 684
$ws.Range("J71").Value = 975.3333
$ws.Range("K71").Value = 6156
$ws.Range("L71").Value = 8777.9997
$ws.Range("M71").Value = -2100
$ws.Range("N71").Value = -16889.9997
# Row 72
$ws.Range("H72").Value = 3834.3333
$ws.Range("I72").Value = 3504.5
$ws.Range("J72").Value = 3999.25
$ws.Range("K72").Value = 31540.5
$ws.Range("L72").Value = 35993.25
$ws.Range("M72").Value = -27484.5
$ws.Range("N72").Value = -44105.25
# Row 122
$ws.Range("H122").Value = 870.9091
$ws.Range("I122").Value = 1046.4286
$ws.Range("J122").Value = 789
$ws.Range("K122").Value = 9417.857399999999
$ws.Range("L122").Value = 7101
$ws.Range("M122").Value = -6967.857399999999
$ws.Range("N122").Value = -12001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4883123
$ws.Range("I70").Value = 6456361
$ws.Range("J70").Value = 6086.8
$ws.Range("K70").Value = 6456361
$ws.Range("L70").Value = 6086.8
$ws.Range("M70").Value = -6456091
$ws.Range("N70").Value = -6626.8
# Row 73
$ws.Range("H73").Value = 4883123
$ws.Range("I73").Value = 6456361
$ws.Range("J73").Value = 6086.8
$ws.Range("K73").Value = 6456361
$ws.Range("L73").Value = 6086.8
$ws.Range("M73").Value = -6455425
$ws.Range("N73").Value = -7958.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4229.5
$ws.Range("I7").Value = 2816.6924
$ws.Range("K7").Value = 2816.6924
$ws.Range("M7").Value = -2704.6924
# Row 46
$ws.Range("H46").Value = 5452.269
$ws.Range("I46").Value = 4726.1665
$ws.Range("J46").Value = 6074.643
$ws.Range("K46").Value = 4726.1665
$ws.Range("L46").Value = 6074.643
$ws.Range("M46").Value = -4538.1665
$ws.Range("N46").Value = -6450.643
# Row 48
$ws.Range("H48").Value = 34254.125
$ws.Range("I48").Value = 32016.5
$ws.Range("K48").Value = 32016.5
$ws.Range("M48").Value = -31355.5
# Row 68
$ws.Range("H68").Value = 2766.818
$ws.Range("I68").Value = 2492.7778
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2492.7778
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1743.7778
$ws.Range("N68").Value = -5498
# Row 71
$ws.Range("H71").Value = 2766.818
$ws.Range("I71").Value = 2492.7778
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 12463.889
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -8719.888999999999
$ws.Range("N71").Value = -27488
# Row 94
$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41352
# Row 126
$ws.Range("H126").Value = 4229.5
$ws.Range("I126").Value = 2816.6924
$ws.Range("K126").Value = 8450.0772
$ws.Range("M126").Value = -5980.0772
# Row 136
$ws.Range("H136").Value = 24856.348
$ws.Range("I136").Value = 28998.578
$ws.Range("J136").Value = 5180.75
$ws.Range("K136").Value = 86995.734
$ws.Range("L136").Value = 15542.25
$ws.Range("M136").Value = -84445.734
$ws.Range("N136").Value = -20642.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
# Row 100
$ws.Range("H100").Value = 1591.75
$ws.Range("I100").Value = 1863.3334
$ws.Range("K100").Value = 3726.6668
$ws.Range("M100").Value = -3185.6668
# Row 136
$ws.Range("H136").Value = 2677.476
$ws.Range("I136").Value = 2347.4
$ws.Range("K136").Value = 7042.200000000001
$ws.Range("M136").Value = -4492.200000000001
